$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-10 from 2023-09-16 (45185) to 2023-10-05 (45204)
foreach ($row in 2..10) {
    $ws.Cells.Item($row, 3).Value = 45204
}
